$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-10 (Sending cluster x Target cluster cross product
# over ECs/FAPs/sCs, with Ligand=Tnc, Receptor=Itgb3 constant).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#          E..T numeric metrics.

$data = @(
    @("ECs","Tnc","Itgb3","ECs",2,0.6666666666666666,1.442371333333333,4.327114,0.03522044016446201,0.03522044016446201,3,1,8.970048,26.910144,0.487108783009476,0.4871087830094759,12.938140093824,116.443260844416,0.01715618574556916,0.01715618574556916),
    @("ECs","Tnc","Itgb3","FAPs",2,0.6666666666666666,1.442371333333333,4.327114,0.03522044016446201,0.03522044016446201,3,1,9.012070666666666,27.036212,0.489390778604016,0.489390778604016,12.99875238357422,116.988771452168,0.01723655863486222,0.01723655863486222),
    @("ECs","Tnc","Itgb3","sCs",2,0.6666666666666666,1.442371333333333,4.327114,0.03522044016446201,0.03522044016446201,3,1,0.4327576666666667,1.298273,0.02350043838650813,0.02350043838650813,0.6241972526802222,5.617775274122,0.0008276957840306357,0.0008276957840306357),
    @("FAPs","Tnc","Itgb3","ECs",3,1,11.331397,33.994191,0.2766948987373093,0.2766948987373092,3,1,8.970048,26.910144,0.487108783009476,0.4871087830094759,101.643174997056,914.788574973504,0.1347805153888609,0.1347805153888609),
    @("FAPs","Tnc","Itgb3","FAPs",3,1,11.331397,33.994191,0.2766948987373093,0.2766948987373092,3,1,9.012070666666666,27.036212,0.489390778604016,0.489390778604016,102.1193505160547,919.074154644492,0.1354119319288112,0.1354119319288111),
    @("FAPs","Tnc","Itgb3","sCs",3,1,11.331397,33.994191,0.2766948987373093,0.2766948987373092,3,1,0.4327576666666667,1.298273,0.02350043838650813,0.02350043838650813,4.903748925793667,44.133740332143,0.006502451419637242,0.006502451419637241),
    @("sCs","Tnc","Itgb3","ECs",3,1,28.17890933333333,84.536728,0.6880846610982287,0.6880846610982286,3,1,8.970048,26.910144,0.487108783009476,0.4871087830094759,252.766169307648,2274.895523768832,0.3351720818750459,0.3351720818750458),
    @("sCs","Tnc","Itgb3","FAPs",3,1,28.17890933333333,84.536728,0.6880846610982287,0.6880846610982286,3,1,9.012070666666666,27.036212,0.489390778604016,0.489390778604016,253.9503222215929,2285.552899994336,0.3367422880403426,0.3367422880403426),
    @("sCs","Tnc","Itgb3","sCs",3,1,28.17890933333333,84.536728,0.6880846610982287,0.6880846610982286,3,1,0.4327576666666667,1.298273,0.02350043838650813,0.02350043838650813,12.19463905230489,109.751751470744,0.01617029118284025,0.01617029118284025)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
